# "add ITI jitter option"
#
# This script:
#  1. Renames two existing nodeHint labels (column C) to re-number the
#     protocol stages, making room for two new EEG stages:
#       "3-时长补充(1k)" -> "4-时长补充(1k)"   (pID 163 rows)
#       "4-阈值阶段(4k)" -> "5-阈值阶段(4k)"   (pID 164 rows)
#  2. Fixes a couple of data values that changed alongside the
#     renumbering (nRepeat for pID161's first row, and the ITI-jitter
#     "f1" window for the loc-1k / dur-1k stages).
#  3. Appends two brand-new "passive EEG" protocol stages (pID 165 =
#     "6-位置阶段EEG(1k)" / SE loc-EEG 1k, and pID 166 =
#     "7-时长补充EEG(1k)" / SE dur-EEG 1k), each following the same
#     11-row ITI-jitter pattern used by the existing stages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Re-label / re-number existing stages (pure text rename, no other
#    column changes for these rows).
# ---------------------------------------------------------------------
$ws.Range("C86:C97").Value = "4-时长补充(1k)"
$ws.Range("C98:C113").Value = "5-阈值阶段(4k)"

# ---------------------------------------------------------------------
# 2) Small data corrections that came along with the renumbering.
# ---------------------------------------------------------------------
# pID 161 ("2-阈值阶段(1k)"): nRepeat for the first row doubled 25 -> 50
$ws.Range("H58").Value = 50

# pID 162 ("3-位置阶段(1k)"): f1 (ITI jitter upper bound) 1030 -> 1015
$ws.Range("L75:L85").Value = 1015

# pID 163 ("4-时长补充(1k)"): f1 (ITI jitter upper bound) 1030 -> 1020
$ws.Range("L87:L97").Value = 1020

# ---------------------------------------------------------------------
# 3) Append the two new passive-EEG stages as rows 114-137.
# ---------------------------------------------------------------------
# Shared "pos" (column N) percentage ladder used by every stage group
# after its first (blank) row.
$posLadder = @(5, 10, 15, 20, 30, 50, 70, 80, 85, 90, 95)

function Add-Stage {
    param($startRow, $pID, $nodeHint, $protocol, $iti, $nRepeat, $durVal)

    for ($i = 0; $i -le 11; $i++) {
        $row = $startRow + $i
        $code = 4 + $i

        $ws.Cells.Item($row, 1).Value = $pID          # A pID
        $ws.Cells.Item($row, 2).Value = "Start-end效应" # B node0Hint
        $ws.Cells.Item($row, 3).Value = $nodeHint       # C nodeHint
        $ws.Cells.Item($row, 4).Value = "passive"       # D apType
        $ws.Cells.Item($row, 5).Value = $protocol       # E protocol
        $ws.Cells.Item($row, 6).Value = $code           # F code
        $ws.Cells.Item($row, 7).Value = $iti            # G ITI
        $ws.Cells.Item($row, 8).Value = $nRepeat        # H nRepeat
        $ws.Cells.Item($row, 9).Value = ""              # I cueLag
        $ws.Cells.Item($row, 10).Value = ""             # J processFcn
        $ws.Cells.Item($row, 11).Value = 1000           # K f0

        if ($i -eq 0) {
            $ws.Cells.Item($row, 12).Value = ""         # L f1
            $ws.Cells.Item($row, 13).Value = ""         # M nChangePeriod
            $ws.Cells.Item($row, 14).Value = ""         # N pos
        }
        else {
            $ws.Cells.Item($row, 12).Value = 1020                 # L f1
            $ws.Cells.Item($row, 13).Value = 20                   # M nChangePeriod
            $ws.Cells.Item($row, 14).Value = $posLadder[$i - 1]   # N pos
        }

        $ws.Cells.Item($row, 15).Value = $durVal        # O dur
    }
}

# pID 165: "6-位置阶段EEG(1k)" / "SE loc-EEG 1k"
Add-Stage 114 165 "6-位置阶段EEG(1k)" "SE loc-EEG 1k" 1.5 40 0.5

# pID 166: "7-时长补充EEG(1k)" / "SE dur-EEG 1k"
Add-Stage 126 166 "7-时长补充EEG(1k)" "SE dur-EEG 1k" 2 40 1
